$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 37.98512966666667
$ws.Range("H2").Value = 113.955389
$ws.Range("I2").Value = 0.5085441461893128
$ws.Range("J2").Value = 0.5085441461893129
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 2928.255539758192
$ws.Range("R2").Value = 26354.29985782373
$ws.Range("S2").Value = 0.1222446897936042
$ws.Range("T2").Value = 0.1222446897936042

# Row 3
$ws.Range("G3").Value = 37.98512966666667
$ws.Range("H3").Value = 113.955389
$ws.Range("I3").Value = 0.5085441461893128
$ws.Range("J3").Value = 0.5085441461893129
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 3858.530889651508
$ws.Range("R3").Value = 34726.77800686357
$ws.Range("S3").Value = 0.1610805154332393
$ws.Range("T3").Value = 0.1610805154332393

# Row 4
$ws.Range("G4").Value = 37.98512966666667
$ws.Range("H4").Value = 113.955389
$ws.Range("I4").Value = 0.5085441461893128
$ws.Range("J4").Value = 0.5085441461893129
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 5394.906008967017
$ws.Range("R4").Value = 48554.15408070316
$ws.Range("S4").Value = 0.2252189409624693
$ws.Range("T4").Value = 0.2252189409624694

# Row 5
$ws.Range("I5").Value = 0.1771904651558058
$ws.Range("J5").Value = 0.1771904651558058
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 1020.283027683631
$ws.Range("R5").Value = 9182.54724915268
$ws.Range("S5").Value = 0.04259333945669377
$ws.Range("T5").Value = 0.04259333945669378

# Row 6
$ws.Range("I6").Value = 0.1771904651558058
$ws.Range("J6").Value = 0.1771904651558058
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.05612478615873692
$ws.Range("T6").Value = 0.05612478615873692

# Row 7
$ws.Range("I7").Value = 0.1771904651558058
$ws.Range("J7").Value = 0.1771904651558058
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.07847233954037508
$ws.Range("T7").Value = 0.0784723395403751

# Row 8
$ws.Range("H8").Value = 70.421093
$ws.Range("I8").Value = 0.3142653886548814
$ws.Range("J8").Value = 0.3142653886548814
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 1809.576163994112
$ws.Range("R8").Value = 16286.18547594701
$ws.Range("S8").Value = 0.07554363812238445
$ws.Range("T8").Value = 0.07554363812238446

# Row 9
$ws.Range("H9").Value = 70.421093
$ws.Range("I9").Value = 0.3142653886548814
$ws.Range("J9").Value = 0.3142653886548814
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.09954304098608344
$ws.Range("T9").Value = 0.09954304098608344

# Row 10
$ws.Range("H10").Value = 70.421093
$ws.Range("I10").Value = 0.3142653886548814
$ws.Range("J10").Value = 0.3142653886548814
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.1391787095464135
$ws.Range("T10").Value = 0.1391787095464135
